$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = -0.089251970589998564
$ws.Range("B1").Value = 0.089114655418271127
$ws.Range("A2").Value = -0.043002458637099394
$ws.Range("B2").Value = 0.042599625591327239
$ws.Range("A3").Value = 0.083884798420214679
$ws.Range("B3").Value = -0.084169029274082874
$ws.Range("A4").Value = -0.19982870760748384
$ws.Range("B4").Value = 0.19877364848234791
$ws.Range("A5").Value = -0.19277364863028534
$ws.Range("B5").Value = 0.19064673046009162
$ws.Range("A6").Value = -0.097522026047394572
$ws.Range("B6").Value = 0.097391070564234905
$ws.Range("A7").Value = -0.077391070743722423
$ws.Range("B7").Value = 0.077084055581940092
$ws.Range("A8").Value = -0.057084055763161778
$ws.Range("B8").Value = 0.0568481813695767
$ws.Range("A9").Value = -0.050848181527304526
$ws.Range("B9").Value = 0.050657755370664859
$ws.Range("A10").Value = -0.044657755530209897
$ws.Range("B10").Value = 0.044631134549334206
$ws.Range("A11").Value = -0.049213311419116224
$ws.Range("B11").Value = 0.049148489689891051
$ws.Range("A12").Value = -0.043148489850041383
$ws.Range("B12").Value = 0.042954381331511193
$ws.Range("A13").Value = -0.036954381493888633
$ws.Range("B13").Value = 0.036902641721296092
$ws.Range("A14").Value = -0.024902641895308442
$ws.Range("B14").Value = 0.02487769319472477
$ws.Range("A15").Value = -0.021053756914719557
$ws.Range("B15").Value = 0.021027982008358315
$ws.Range("A16").Value = -0.015027982172391763
$ws.Range("B16").Value = 0.015004746394037927
$ws.Range("A17").Value = -0.0090047465587250741
$ws.Range("B17").Value = 0.0089999998296708128
$ws.Range("A18").Value = -0.036112196934027452
$ws.Range("B18").Value = 0.036097023388052207
$ws.Range("A19").Value = -0.027097023539131904
$ws.Range("B19").Value = 0.027014138661656695
$ws.Range("A20").Value = -0.018014138813995828
$ws.Range("B20").Value = 0.018004317394058944
$ws.Range("A21").Value = -0.0090043175465623904
$ws.Range("B21").Value = 0.008999999847353557
$ws.Range("A22").Value = -0.11748442416330107
$ws.Range("B22").Value = 0.11706984637147322
$ws.Range("A23").Value = -0.084625995619954431
$ws.Range("B23").Value = 0.084125246721977298
$ws.Range("A24").Value = -0.042125246937942507
$ws.Range("B24").Value = 0.041999999782953878
$ws.Range("A25").Value = -0.087124704565987088
$ws.Range("B25").Value = 0.087019074883723846
$ws.Range("A26").Value = -0.081019075037449539
$ws.Range("B26").Value = 0.080890911348472372
$ws.Range("A27").Value = -0.074890911502912161
$ws.Range("B27").Value = 0.074481302412336525
$ws.Range("A28").Value = -0.068481302569525226
$ws.Range("B28").Value = 0.068219463514950895
$ws.Range("A29").Value = -0.056219463684966442
$ws.Range("B29").Value = 0.056110580580451597
$ws.Range("A30").Value = -0.03611058076598761
$ws.Range("B30").Value = 0.035845768871807682
$ws.Range("A31").Value = -0.020845769049397944
$ws.Range("B31").Value = 0.020763062680261513
$ws.Range("A32").Value = -0.006001016956275862
$ws.Range("B32").Value = 0.0059999998386244258

# Set column B width (target OOXML stored width 14.7109375 -> achievable quantized width 14.666666666666666 via ColumnWidth=13.8)
$ws.Columns.Item(2).ColumnWidth = 13.8

Write-Output "Edit applied"
